# Auto-generated edit script applying the Cactuar_Profits.xlsx numeric update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1264.25
$ws.Range("I15").Value = 1264.25
$ws.Range("K15").Value = 3792.75
$ws.Range("M15").Value = -3623.75
$ws.Range("H32").Value = 62505460
$ws.Range("I32").Value = 250002240
$ws.Range("J32").Value = 6533.5
$ws.Range("K32").Value = 250002240
$ws.Range("L32").Value = 6533.5
$ws.Range("M32").Value = -250001914
$ws.Range("N32").Value = -7185.5
$ws.Range("H33").Value = 534.3333
$ws.Range("I33").Value = 226.375
$ws.Range("K33").Value = 226.375
$ws.Range("M33").Value = 2.625
$ws.Range("H40").Value = 2160362.5
$ws.Range("J40").Value = 3518223
$ws.Range("L40").Value = 3518223
$ws.Range("N40").Value = -3518573
$ws.Range("H51").Value = 7226.4546
$ws.Range("I51").Value = 4562.8
$ws.Range("J51").Value = 9446.166999999999
$ws.Range("K51").Value = 4562.8
$ws.Range("L51").Value = 9446.166999999999
$ws.Range("M51").Value = -4078.8
$ws.Range("N51").Value = -10414.167
$ws.Range("H70").Value = 4398.278
$ws.Range("I70").Value = 2874.3333
$ws.Range("J70").Value = 5922.222
$ws.Range("K70").Value = 8622.999899999999
$ws.Range("L70").Value = 17766.666
$ws.Range("M70").Value = -8352.999899999999
$ws.Range("N70").Value = -18306.666
$ws.Range("H73").Value = 4398.278
$ws.Range("I73").Value = 2874.3333
$ws.Range("J73").Value = 5922.222
$ws.Range("K73").Value = 8622.999899999999
$ws.Range("L73").Value = 17766.666
$ws.Range("M73").Value = -7686.999899999999
$ws.Range("N73").Value = -19638.666
$ws.Range("H80").Value = 31250812
$ws.Range("J80").Value = 1145.6
$ws.Range("L80").Value = 3436.8
$ws.Range("N80").Value = -5432.799999999999
$ws.Range("H83").Value = 31250812
$ws.Range("J83").Value = 1145.6
$ws.Range("L83").Value = 10310.4
$ws.Range("N83").Value = -20294.4
$ws.Range("H132").Value = 119074.875
$ws.Range("I132").Value = 162896.14
$ws.Range("K132").Value = 488688.42
$ws.Range("M132").Value = -486158.42
$ws.Range("H135").Value = 4433.1387
$ws.Range("I135").Value = 1238.6842
$ws.Range("K135").Value = 11148.1578
$ws.Range("M135").Value = -8613.157799999999
$ws.Range("H141").Value = 3990
$ws.Range("I141").Value = 4034.88
$ws.Range("K141").Value = 12104.64
$ws.Range("M141").Value = -6924.639999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1343688.2
$ws.Range("I2").Value = 1940442.9
$ws.Range("K2").Value = 1940442.9
$ws.Range("M2").Value = -1940329.9
$ws.Range("H32").Value = 16126.28
$ws.Range("I32").Value = 15554
$ws.Range("K32").Value = 15554
$ws.Range("M32").Value = -15267
$ws.Range("H61").Value = 5798.75
$ws.Range("I61").Value = 5620.2974
$ws.Range("K61").Value = 5620.2974
$ws.Range("M61").Value = -5408.2974
$ws.Range("H74").Value = 16668380
$ws.Range("J74").Value = 2054
$ws.Range("L74").Value = 2054
$ws.Range("N74").Value = -3802
$ws.Range("H77").Value = 16668380
$ws.Range("J77").Value = 2054
$ws.Range("L77").Value = 10270
$ws.Range("N77").Value = -19006
$ws.Range("H116").Value = 1343688.2
$ws.Range("I116").Value = 1940442.9
$ws.Range("K116").Value = 1940442.9
$ws.Range("M116").Value = -1938148.9
$ws.Range("H136").Value = 5798.75
$ws.Range("I136").Value = 5620.2974
$ws.Range("K136").Value = 16860.8922
$ws.Range("M136").Value = -14310.8922

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1343688.2
$ws.Range("I3").Value = 1940442.9
$ws.Range("K3").Value = 1940442.9
$ws.Range("M3").Value = -1940328.9
$ws.Range("H20").Value = 4589.4585
$ws.Range("I20").Value = 4228.25
$ws.Range("J20").Value = 5311.875
$ws.Range("K20").Value = 4228.25
$ws.Range("L20").Value = 5311.875
$ws.Range("M20").Value = -3981.25
$ws.Range("N20").Value = -5805.875
$ws.Range("H86").Value = 1796.4828
$ws.Range("I86").Value = 1451.6364
$ws.Range("K86").Value = 1451.6364
$ws.Range("M86").Value = -328.6364000000001
$ws.Range("H89").Value = 1796.4828
$ws.Range("I89").Value = 1451.6364
$ws.Range("K89").Value = 7258.182000000001
$ws.Range("M89").Value = -1642.182000000001
$ws.Range("H105").Value = 1402.9286
$ws.Range("I105").Value = 1402.9286
$ws.Range("K105").Value = 1402.9286
$ws.Range("M105").Value = 344.0714
$ws.Range("H107").Value = 1309.0526
$ws.Range("I107").Value = 1319.3846
$ws.Range("J107").Value = 1286.6666
$ws.Range("K107").Value = 1319.3846
$ws.Range("L107").Value = 1286.6666
$ws.Range("M107").Value = 600.6153999999999
$ws.Range("N107").Value = -5126.6666
$ws.Range("H134").Value = 1223.2222
$ws.Range("I134").Value = 1252.1923
$ws.Range("K134").Value = 3756.5769
$ws.Range("M134").Value = -1221.5769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19236686
$ws.Range("I31").Value = 76925170
$ws.Range("K31").Value = 76925170
$ws.Range("M31").Value = -76924875
$ws.Range("H34").Value = 19236686
$ws.Range("I34").Value = 76925170
$ws.Range("K34").Value = 76925170
$ws.Range("M34").Value = -76924968
$ws.Range("H105").Value = 1749638.9
$ws.Range("I105").Value = 2067391.4
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 2067391.4
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -2065644.4
$ws.Range("N105").Value = -5494
$ws.Range("H134").Value = 1716.9333
$ws.Range("I134").Value = 1696.8214
$ws.Range("K134").Value = 5090.4642
$ws.Range("M134").Value = -2555.4642

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1239.6666
$ws.Range("J98").Value = 1272.75
$ws.Range("L98").Value = 3818.25
$ws.Range("N98").Value = -6814.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2038
$ws.Range("I113").Value = 1792.75
$ws.Range("K113").Value = 1792.75
$ws.Range("M113").Value = 377.25
$ws.Range("H134").Value = 5000000
$ws.Range("J134").Value = 5000000
$ws.Range("L134").Value = 15000000
$ws.Range("N134").Value = -15005070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5085.125
$ws.Range("I7").Value = 3431.3
$ws.Range("K7").Value = 3431.3
$ws.Range("M7").Value = -3319.3
$ws.Range("H22").Value = 1048.8572
$ws.Range("I22").Value = 914
$ws.Range("J22").Value = 1183.7142
$ws.Range("K22").Value = 914
$ws.Range("L22").Value = 1183.7142
$ws.Range("M22").Value = -619
$ws.Range("N22").Value = -1773.7142
$ws.Range("H27").Value = 1048.8572
$ws.Range("I27").Value = 914
$ws.Range("J27").Value = 1183.7142
$ws.Range("K27").Value = 914
$ws.Range("L27").Value = 1183.7142
$ws.Range("M27").Value = -807
$ws.Range("N27").Value = -1397.7142
$ws.Range("H40").Value = 44449444
$ws.Range("J40").Value = 83339336
$ws.Range("L40").Value = 83339336
$ws.Range("N40").Value = -83339608
$ws.Range("H46").Value = 6777
$ws.Range("J46").Value = 7045.143
$ws.Range("L46").Value = 7045.143
$ws.Range("N46").Value = -7421.143
$ws.Range("H55").Value = 327.4138
$ws.Range("I55").Value = 179.46666
$ws.Range("K55").Value = 179.46666
$ws.Range("M55").Value = -6.46665999999999
$ws.Range("H82").Value = 2842443.8
$ws.Range("J82").Value = 1765.2222
$ws.Range("L82").Value = 1765.2222
$ws.Range("N82").Value = -2487.2222
$ws.Range("H85").Value = 2842443.8
$ws.Range("J85").Value = 1765.2222
$ws.Range("L85").Value = 1765.2222
$ws.Range("N85").Value = -4261.2222
$ws.Range("H108").Value = 29726
$ws.Range("J108").Value = 29726
$ws.Range("L108").Value = 29726
$ws.Range("N108").Value = -37406
$ws.Range("H122").Value = 4583.069
$ws.Range("I122").Value = 3626.238
$ws.Range("K122").Value = 10878.714
$ws.Range("M122").Value = -8428.714
$ws.Range("H126").Value = 5085.125
$ws.Range("I126").Value = 3431.3
$ws.Range("K126").Value = 10293.9
$ws.Range("M126").Value = -7823.900000000001
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2982973.8
$ws.Range("I81").Value = 4171564.5
$ws.Range("K81").Value = 8343129
$ws.Range("M81").Value = -8342068
$ws.Range("H84").Value = 2982973.8
$ws.Range("I84").Value = 4171564.5
$ws.Range("K84").Value = 41715645
$ws.Range("M84").Value = -41710341
$ws.Range("H123").Value = 76747.25
$ws.Range("J123").Value = 76747.25
$ws.Range("L123").Value = 76747.25
$ws.Range("N123").Value = -86547.25
$ws.Range("H132").Value = 12501145
$ws.Range("I132").Value = 1037.8235
$ws.Range("K132").Value = 3113.4705
$ws.Range("M132").Value = -583.4704999999999
$ws.Range("H136").Value = 5690.63
$ws.Range("I136").Value = 2264.0715
$ws.Range("J136").Value = 10051.704
$ws.Range("K136").Value = 6792.2145
$ws.Range("L136").Value = 30155.112
$ws.Range("M136").Value = -4242.2145

